$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "Lattitude" -> "Latitude"
$ws.Range("D1").Value = "Latitude"

# Fix city names (remove stray trailing space/nbsp, remove comma)
$ws.Range("A8").Value = "Den Haag"
$ws.Range("A15").Value = "Schiphol Haarlemmermeer"

# Replace the one-hot header row (F1:Z1), which held city-name strings,
# with plain sequential numbers 0-20
$ws.Range("F1").Value = 0
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = 2
$ws.Range("I1").Value = 3
$ws.Range("J1").Value = 4
$ws.Range("K1").Value = 5
$ws.Range("L1").Value = 6
$ws.Range("M1").Value = 7
$ws.Range("N1").Value = 8
$ws.Range("O1").Value = 9
$ws.Range("P1").Value = 10
$ws.Range("Q1").Value = 11
$ws.Range("R1").Value = 12
$ws.Range("S1").Value = 13
$ws.Range("T1").Value = 14
$ws.Range("U1").Value = 15
$ws.Range("V1").Value = 16
$ws.Range("W1").Value = 17
$ws.Range("X1").Value = 18
$ws.Range("Y1").Value = 19
$ws.Range("Z1").Value = 20

# Move selection to A15 (matches author's recorded cursor position)
$ws.Range("A15").Select() | Out-Null
